$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (artifact_type), shifting
# artifact_type/text right to make room for the new tomb_code column.
$ws.Columns.Item(2).Insert()

# Header for the new column
$ws.Cells.Item(1, 2).Value = "tomb_code"

# Tablet rows (2-9) get sequential tomb codes 1-8
$ws.Cells.Item(2, 2).Value = 1
$ws.Cells.Item(3, 2).Value = 2
$ws.Cells.Item(4, 2).Value = 3
$ws.Cells.Item(5, 2).Value = 4
$ws.Cells.Item(6, 2).Value = 5
$ws.Cells.Item(7, 2).Value = 6
$ws.Cells.Item(8, 2).Value = 7
$ws.Cells.Item(9, 2).Value = 8

# Relic rows (10-17) have no tomb code -- leave column B blank there.

# Column widths: A and C keep their auto-fit ("best fit") widths; the
# newly inserted column B gets an explicit (non-autofit) width that
# matches column A, mirroring the author's manual column resize.
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(1).ColumnWidth

# Adjust selection to match target workbook state
$ws.Range("B4").Select()
